# Update cryptocurrency price / 1h-volume figures scraped from coinranking.com,
# plus re-sort two pairs of adjacent rows whose ranking flipped between runs
# (Polkadot/Polygon at rows 13-14, ShibaInu/BitcoinCash at rows 18-19,
# FraxShare/Aave at rows 43-44) and replace SynthetixNetwork (row 51) with Aptos.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (Coin) and C (Link) are plain text already; column E (Volume) is
# always padded with spaces and a "%" sign so it never parses as a number.
# Column D (Price) sometimes contains a string that Excel would otherwise
# auto-convert to a floating-point Double (losing the original text layout),
# so those cells are pre-formatted as Text before the value is written.

$ws.Range('D2').Value = '29.619.26'
$ws.Range('E2').Value = '  +1.30%  '
$ws.Range('D3').Value = '1.882.81'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9990'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7273'
$ws.Range('E5').Value = '  +3.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '239.58'
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07859'
$ws.Range('E8').Value = '  -3.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3095'
$ws.Range('E9').Value = '  +2.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.29'
$ws.Range('E10').Value = '  +8.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08229'
$ws.Range('E11').Value = '  +0.98%  '
$ws.Range('D12').Value = '1.887.67'
$ws.Range('E12').Value = '  +1.62%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7295'
$ws.Range('E13').Value = '  +3.39%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.286'
$ws.Range('E14').Value = '  +2.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '90.27'
$ws.Range('E15').Value = '  +1.31%  '
$ws.Range('D16').Value = '29.556.73'
$ws.Range('E16').Value = '  +1.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.885'
$ws.Range('E17').Value = '  +1.95%  '
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '243.17'
$ws.Range('E18').Value = '  +3.16%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007883'
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.43'
$ws.Range('E20').Value = '  +0.73%  '
$ws.Range('D21').Value = '2.115.61'
$ws.Range('E21').Value = '  +0.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9988'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9988'
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.789'
$ws.Range('E24').Value = '  +5.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1505'
$ws.Range('E25').Value = '  +4.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.94'
$ws.Range('E26').Value = '  +1.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.008'
$ws.Range('E27').Value = '  +0.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.34'
$ws.Range('E28').Value = '  +1.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.963'
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.367'
$ws.Range('E30').Value = '  -4.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.484'
$ws.Range('E31').Value = '  +0.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.371'
$ws.Range('E32').Value = '  -0.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.111'
$ws.Range('E33').Value = '  +1.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05266'
$ws.Range('E34').Value = '  +1.54%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.202'
$ws.Range('E35').Value = '  +2.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7187'
$ws.Range('E36').Value = '  +1.85%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.004'
$ws.Range('E37').Value = '  +0.62%  '
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01868'
$ws.Range('E39').Value = '  +1.48%  '
$ws.Range('E40').Value = '  -0.55%  '
$ws.Range('D41').Value = '1.178.16'
$ws.Range('E41').Value = '  +4.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9140'
$ws.Range('E42').Value = '  -0.56%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '72.27'
$ws.Range('E43').Value = '  +3.09%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.007'
$ws.Range('E44').Value = '  +2.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4328'
$ws.Range('E45').Value = '  +1.38%  '
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.69'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5345'
$ws.Range('E48').Value = '  -1.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.782'
$ws.Range('E49').Value = '  +1.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.257'
$ws.Range('E50').Value = '  +1.00%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.109'
$ws.Range('E51').Value = '  +2.44%  '
